$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new "Gender" column before the existing "Email" column (old D,
# which shifts to E). This pushes Lada/Phone/Cellphone/BornDate/NationalId
# one column to the right as well (old E..I -> new F..J).
# ---------------------------------------------------------------------------
$ws.Columns("D:D").Insert()

# Preserve the original hyperlink-cell formatting (distinct styles used on
# row 2 vs rows 3-5) before Hyperlinks.Add() stomps it with the built-in
# "Hyperlink" style. Stash copies of the current formats in a scratch cell
# far away from the used range, then paste them back after rewiring the
# hyperlinks, finally clearing the scratch cell so it doesn't linger.
$ws.Range("E2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("Z2").PasteSpecial(-4122)

# The column insert shifted the underlying cells (and the mailto hyperlink
# relationships) from D to E automatically, but the <hyperlink ref="..."/>
# entries themselves still point at the old D column, so rebuild them
# against column E.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:enriquecarrillo119999@gmail.com", "", "", "enriquecarrillo119999@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:edgarcarrillo119999@gmail.com", "", "", "edgarcarrillo119999@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:myriamcarrillo119999@gmail.com", "", "", "myriamcarrillo119999@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:juliocarrillo119999@gmail.com", "", "", "juliocarrillo119999@gmail.com")

# Restore the pre-hyperlink formatting.
$ws.Range("Z1").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("E5").PasteSpecial(-4122)

$ws.Range("Z1:Z2").Clear()

# New Gender column is much narrower than the old wide Email column.
$ws.Columns("D:D").ColumnWidth = 6.95

# Re-touch the NationalId column (shifted from I to J) so its custom width
# survives the save round-trip.
$ws.Columns("J:J").ColumnWidth = 9.49

# Header + data for the new Gender column.
$ws.Range("D1").Value = "Gender"
$ws.Range("D2").Value = "M"
$ws.Range("D3").Value = "M"
$ws.Range("D4").Value = "F"
$ws.Range("D5").Value = "M"

# Match the author's final cursor position.
$ws.Range("D7").Select()
